# Generate Report for Handback
# Updates the "generate date" timestamps that are refreshed each time the
# handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-08-29 21:19:39"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$wsZhCn.Range("H2").Value = "2016-08-29 21:19:35"
$wsZhCn.Range("K2").Value = "2016-08-29 21:19:53"

# de-de sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$wsDeDe.Range("H2").Value = "2016-08-29 21:19:39"
$wsDeDe.Range("K2").Value = "2016-08-29 21:20:02"
